{"js": "// Replace the 25 \"NNN\u00f7N=NN, N\" answer strings in the 5x5 practice table.\n// The table has 20 rows total: a data row followed by 3 blank rows,\n// repeated 5 times (data rows at index 0, 4, 8, 12, 16), 5 columns each.\n// Each (row, col) pair below gives the expected current text (oldText)\n// and the text it must become (newText); oldText is used purely as a\n// sanity check so the correct cell is never ambiguous.\nconst replacements = [\n  { row: 0, col: 0, oldText: \"675\u00f79=75, 0\", newText: \"363\u00f72=181, 1\" },\n  { row: 0, col: 1, oldText: \"422\u00f73=140, 2\", newText: \"369\u00f77=52, 5\" },\n  { row: 0, col: 2, oldText: \"815\u00f73=271, 2\", newText: \"134\u00f72=67, 0\" },\n  { row: 0, col: 3, oldText: \"482\u00f77=68, 6\", newText: \"832\u00f77=118, 6\" },\n  { row: 0, col: 4, oldText: \"132\u00f78=16, 4\", newText: \"659\u00f74=164, 3\" },\n  { row: 4, col: 0, oldText: \"719\u00f72=359, 1\", newText: \"348\u00f79=38, 6\" },\n  { row: 4, col: 1, oldText: \"520\u00f76=86, 4\", newText: \"654\u00f76=109, 0\" },\n  { row: 4, col: 2, oldText: \"399\u00f72=199, 1\", newText: \"929\u00f73=309, 2\" },\n  { row: 4, col: 3, oldText: \"114\u00f73=38, 0\", newText: \"674\u00f78=84, 2\" },\n  { row: 4, col: 4, oldText: \"570\u00f79=63, 3\", newText: \"230\u00f77=32, 6\" },\n  { row: 8, col: 0, oldText: \"861\u00f79=95, 6\", newText: \"875\u00f72=437, 1\" },\n  { row: 8, col: 1, oldText: \"595\u00f73=198, 1\", newText: \"268\u00f77=38, 2\" },\n  { row: 8, col: 2, oldText: \"519\u00f75=103, 4\", newText: \"231\u00f78=28, 7\" },\n  { row: 8, col: 3, oldText: \"219\u00f74=54, 3\", newText: \"345\u00f78=43, 1\" },\n  { row: 8, col: 4, oldText: \"609\u00f79=67, 6\", newText: \"355\u00f76=59, 1\" },\n  { row: 12, col: 0, oldText: \"866\u00f74=216, 2\", newText: \"499\u00f74=124, 3\" },\n  { row: 12, col: 1, oldText: \"350\u00f74=87, 2\", newText: \"490\u00f73=163, 1\" },\n  { row: 12, col: 2, oldText: \"823\u00f77=117, 4\", newText: \"120\u00f76=20, 0\" },\n  { row: 12, col: 3, oldText: \"149\u00f73=49, 2\", newText: \"598\u00f78=74, 6\" },\n  { row: 12, col: 4, oldText: \"194\u00f79=21, 5\", newText: \"604\u00f77=86, 2\" },\n  { row: 16, col: 0, oldText: \"719\u00f73=239, 2\", newText: \"848\u00f77=121, 1\" },\n  { row: 16, col: 1, oldText: \"922\u00f78=115, 2\", newText: \"760\u00f76=126, 4\" },\n  { row: 16, col: 2, oldText: \"963\u00f77=137, 4\", newText: \"864\u00f74=216, 0\" },\n  { row: 16, col: 3, oldText: \"388\u00f79=43, 1\", newText: \"109\u00f76=18, 1\" },\n  { row: 16, col: 4, oldText: \"231\u00f78=28, 7\", newText: \"562\u00f76=93, 4\" },\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\nconst table = tables.items[0];\n\n// Load the current text of every target cell up front so we can verify\n// we are editing the right place before writing anything.\nconst cells = replacements.map((r) => table.getCell(r.row, r.col));\ncells.forEach((cell) => cell.body.load(\"text\"));\nawait context.sync();\n\nfor (let i = 0; i < replacements.length; i++) {\n  const { oldText, newText } = replacements[i];\n  const actual = cells[i].body.text.trim();\n  if (actual !== oldText) {\n    throw new Error(\n      `Unexpected cell text at index ${i}: expected \"${oldText}\" got \"${actual}\"`\n    );\n  }\n  cells[i].value = newText;\n}\n\nawait context.sync();\n", "ps1": "# Replace the 25 \"NNN\u00f7N=NN, N\" answer strings in the 5x5 practice table.\n# The table has 20 rows total: a data row followed by 3 blank rows,\n# repeated 5 times (data rows at Word's 1-based row 1, 5, 9, 13, 17),\n# 5 columns each. We address cells by (row, col) position rather than a\n# document-wide Find/Replace because a couple of the target strings\n# duplicate other rows' source strings (\"231\u00f78=28, 7\" is both the OLD\n# text of one cell and the NEW text of a later cell) - a global replace\n# could touch the wrong cell once the first write lands.\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$replacements = @(\n    @{ Row = 1; Col = 1; OldText = \"675\u00f79=75, 0\"; NewText = \"363\u00f72=181, 1\" },\n    @{ Row = 1; Col = 2; OldText = \"422\u00f73=140, 2\"; NewText = \"369\u00f77=52, 5\" },\n    @{ Row = 1; Col = 3; OldText = \"815\u00f73=271, 2\"; NewText = \"134\u00f72=67, 0\" },\n    @{ Row = 1; Col = 4; OldText = \"482\u00f77=68, 6\"; NewText = \"832\u00f77=118, 6\" },\n    @{ Row = 1; Col = 5; OldText = \"132\u00f78=16, 4\"; NewText = \"659\u00f74=164, 3\" },\n    @{ Row = 5; Col = 1; OldText = \"719\u00f72=359, 1\"; NewText = \"348\u00f79=38, 6\" },\n    @{ Row = 5; Col = 2; OldText = \"520\u00f76=86, 4\"; NewText = \"654\u00f76=109, 0\" },\n    @{ Row = 5; Col = 3; OldText = \"399\u00f72=199, 1\"; NewText = \"929\u00f73=309, 2\" },\n    @{ Row = 5; Col = 4; OldText = \"114\u00f73=38, 0\"; NewText = \"674\u00f78=84, 2\" },\n    @{ Row = 5; Col = 5; OldText = \"570\u00f79=63, 3\"; NewText = \"230\u00f77=32, 6\" },\n    @{ Row = 9; Col = 1; OldText = \"861\u00f79=95, 6\"; NewText = \"875\u00f72=437, 1\" },\n    @{ Row = 9; Col = 2; OldText = \"595\u00f73=198, 1\"; NewText = \"268\u00f77=38, 2\" },\n    @{ Row = 9; Col = 3; OldText = \"519\u00f75=103, 4\"; NewText = \"231\u00f78=28, 7\" },\n    @{ Row = 9; Col = 4; OldText = \"219\u00f74=54, 3\"; NewText = \"345\u00f78=43, 1\" },\n    @{ Row = 9; Col = 5; OldText = \"609\u00f79=67, 6\"; NewText = \"355\u00f76=59, 1\" },\n    @{ Row = 13; Col = 1; OldText = \"866\u00f74=216, 2\"; NewText = \"499\u00f74=124, 3\" },\n    @{ Row = 13; Col = 2; OldText = \"350\u00f74=87, 2\"; NewText = \"490\u00f73=163, 1\" },\n    @{ Row = 13; Col = 3; OldText = \"823\u00f77=117, 4\"; NewText = \"120\u00f76=20, 0\" },\n    @{ Row = 13; Col = 4; OldText = \"149\u00f73=49, 2\"; NewText = \"598\u00f78=74, 6\" },\n    @{ Row = 13; Col = 5; OldText = \"194\u00f79=21, 5\"; NewText = \"604\u00f77=86, 2\" },\n    @{ Row = 17; Col = 1; OldText = \"719\u00f73=239, 2\"; NewText = \"848\u00f77=121, 1\" },\n    @{ Row = 17; Col = 2; OldText = \"922\u00f78=115, 2\"; NewText = \"760\u00f76=126, 4\" },\n    @{ Row = 17; Col = 3; OldText = \"963\u00f77=137, 4\"; NewText = \"864\u00f74=216, 0\" },\n    @{ Row = 17; Col = 4; OldText = \"388\u00f79=43, 1\"; NewText = \"109\u00f76=18, 1\" },\n    @{ Row = 17; Col = 5; OldText = \"231\u00f78=28, 7\"; NewText = \"562\u00f76=93, 4\" }\n)\n\nforeach ($item in $replacements) {\n    $cell = $t.Cell($item.Row, $item.Col)\n    $r = $cell.Range\n    # Trim the trailing cell-end mark(s) so we compare/replace only the\n    # visible text, not Word's end-of-cell/end-of-row control characters.\n    $r.End = $r.End - 1\n    $actual = $r.Text\n    if ($actual -ne $item.OldText) {\n        throw \"Unexpected cell text at row $($item.Row) col $($item.Col): expected '$($item.OldText)' got '$actual'\"\n    }\n    $r.Text = $item.NewText\n}\n"}
